$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Avaliações") holds Brazilian-formatted integers (e.g. "1.293")
# that use "." as a thousands separator, not a decimal point. Force the
# cells to Text format before writing so Excel stores them verbatim instead
# of re-parsing them as numbers.

# Row 2
$ws.Cells.Item(2, 5).Value = "Rua Pedro Doll, 564"

# Row 3
$ws.Cells.Item(3, 3).Value = "5,0"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "151"
$ws.Cells.Item(3, 5).Value = "Rua Francisca Júlia, 524"

# Row 4
$ws.Cells.Item(4, 3).Value = "5,0"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "118"
$ws.Cells.Item(4, 5).Value = "R. Conselheiro Saraiva, 207"

# Row 5
$ws.Cells.Item(5, 1).Value = "Black Fitness Club"
$ws.Cells.Item(5, 3).Value = "4,7"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "196"
$ws.Cells.Item(5, 5).Value = "R. Gaspar Soares, 88"

# Row 6
$ws.Cells.Item(6, 1).Value = "BlueFit Santana"
$ws.Cells.Item(6, 3).Value = "4,5"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "691"
$ws.Cells.Item(6, 5).Value = "R. Dr. Zuquim, 1872"

# Row 7
$ws.Cells.Item(7, 1).Value = "Health Academy"
$ws.Cells.Item(7, 3).Value = "5,0"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "18"
$ws.Cells.Item(7, 5).Value = "R. Voluntários da Pátria, 2468 - Conjunto 167"

# Row 8
$ws.Cells.Item(8, 1).Value = "Tracer Parkour - Santana"
$ws.Cells.Item(8, 3).Value = "4,8"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "52"
$ws.Cells.Item(8, 5).Value = "R. Dr. Zuquim, 859"

# Row 9
$ws.Cells.Item(9, 1).Value = "Silver | Fitness & Spa | Santana - Zona Norte - SP"
$ws.Cells.Item(9, 3).Value = "4,8"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "110"
$ws.Cells.Item(9, 5).Value = "R. Conselheiro Moreira de Barros, 636"

# Row 10
$ws.Cells.Item(10, 1).Value = "Academia Butterfly"
$ws.Cells.Item(10, 3).Value = "4,6"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "148"
$ws.Cells.Item(10, 5).Value = "R. Conselheiro Moreira de Barros, 221"

# Row 11
$ws.Cells.Item(11, 3).Value = "3,9"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "156"
$ws.Cells.Item(11, 5).Value = "Rua Dr. Olavo Egídio, 632"

# Row 12
$ws.Cells.Item(12, 3).Value = "4,6"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "18"
$ws.Cells.Item(12, 5).Value = "Rua Benvinda Aparecida de Abreu Leme, 43 - Sala 3"

# Row 13
$ws.Cells.Item(13, 3).Value = "4,5"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "233"
$ws.Cells.Item(13, 5).Value = "R. Ten. Blum, 93"

# Row 14
$ws.Cells.Item(14, 3).Value = "5,0"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "151"
$ws.Cells.Item(14, 5).Value = "R. Francisca Júlia, 524"

# Row 15
$ws.Cells.Item(15, 3).Value = "4,7"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "510"
$ws.Cells.Item(15, 5).Value = "Rua Chemin Del Pra, 58"

# Row 16
$ws.Cells.Item(16, 3).Value = "4,9"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "166"
$ws.Cells.Item(16, 5).Value = "R. Conselheiro Moreira de Barros, 711"

# Row 17
$ws.Cells.Item(17, 3).Value = "4,8"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "62"
$ws.Cells.Item(17, 5).Value = "Av. Nova Cantareira, 624"

# Row 18
$ws.Cells.Item(18, 3).Value = "4,1"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "598"
$ws.Cells.Item(18, 5).Value = "Av. Luiz Dumont Villares, 200"

# Row 19
$ws.Cells.Item(19, 3).Value = "4,9"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "1.293"
$ws.Cells.Item(19, 5).Value = "R. Mateus Leme, 114"

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "254"
$ws.Cells.Item(20, 5).Value = "R. Voluntários da Pátria, 1884"

# Row 21
$ws.Cells.Item(21, 5).Value = "R. Pedro Doll, 564"

# Row 22
$ws.Cells.Item(22, 5).Value = "R. Antônio Pereira de Sousa, 227"

# Row 23
$ws.Cells.Item(23, 1).Value = "Angusteria"
$ws.Cells.Item(23, 3).Value = "4,7"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "443"
$ws.Cells.Item(23, 5).Value = "R. Duarte de Azevedo, 468"

# Row 24
$ws.Cells.Item(24, 1).Value = "Lassù"
$ws.Cells.Item(24, 3).Value = "4,6"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "3.332"
$ws.Cells.Item(24, 5).Value = "R. Conselheiro Saraiva, 207"

# Row 25
$ws.Cells.Item(25, 1).Value = "Famiglia Mancini Trattoria"
$ws.Cells.Item(25, 3).Value = "4,7"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "27.903"
$ws.Cells.Item(25, 5).Value = "Rua Avanhandava, 81, Bela Vista"

# Row 26
$ws.Cells.Item(26, 1).Value = "La Mordida"
$ws.Cells.Item(26, 3).Value = "4,6"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "618"
$ws.Cells.Item(26, 5).Value = "R. Prof. Lourival Gomes Machado, 274 - 272"

# Row 27
$ws.Cells.Item(27, 1).Value = "Rainha da Jovita"
$ws.Cells.Item(27, 3).Value = "4,6"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "172"
$ws.Cells.Item(27, 5).Value = "Rua Dr. Olavo Egídio, 449"

# Row 28
$ws.Cells.Item(28, 1).Value = "La Ficazza"
$ws.Cells.Item(28, 3).Value = "4,9"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "730"
$ws.Cells.Item(28, 5).Value = "R. Dr. César, 1160"

# Row 29
$ws.Cells.Item(29, 1).Value = "Sabor de Santana"
$ws.Cells.Item(29, 3).Value = "4,6"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "194"
$ws.Cells.Item(29, 5).Value = "R. Duarte de Azevedo, 352"

# Row 30
$ws.Cells.Item(30, 1).Value = "LeLui Bar e Cozinha"
$ws.Cells.Item(30, 3).Value = "4,8"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.018"
$ws.Cells.Item(30, 5).Value = "R. Jacuna, 302"

# Row 31
$ws.Cells.Item(31, 1).Value = "Dom ramiro"
$ws.Cells.Item(31, 3).Value = "4,5"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "3.579"
$ws.Cells.Item(31, 5).Value = "R. Dr. César, 105"

# Row 32
$ws.Cells.Item(32, 3).Value = "4,3"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "114"
$ws.Cells.Item(32, 5).Value = "R. Dr. Zuquim, 1941"

# Row 33
$ws.Cells.Item(33, 1).Value = "Salvi Café e Cozinha - Santana"
$ws.Cells.Item(33, 3).Value = "4,6"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.216"
$ws.Cells.Item(33, 5).Value = "R. Tupiguaes, 140"

# Row 34
$ws.Cells.Item(34, 1).Value = "Pecorino Cucina Mediterranea Braz Leme: Restaurante, Delivery São Paulo SP"
$ws.Cells.Item(34, 3).Value = "4,6"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.549"
$ws.Cells.Item(34, 5).Value = "Av. Braz Leme, 1200"

# Row 35
$ws.Cells.Item(35, 1).Value = "Santo Mar Restaurante — Unidade Santana"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.898"
$ws.Cells.Item(35, 5).Value = "Av. Luiz Dumont Villares, 1306"

# Row 36
$ws.Cells.Item(36, 3).Value = "4,2"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "363"
$ws.Cells.Item(36, 5).Value = "R. Prof. Lourival Gomes Machado, 273"

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.741"
$ws.Cells.Item(37, 5).Value = "R. Lucas de Freitas de Azevedo, 87"

# Row 38
$ws.Cells.Item(38, 1).Value = "Quintal Brasil"
$ws.Cells.Item(38, 3).Value = "4,2"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "770"
$ws.Cells.Item(38, 5).Value = "R. Alfredo Pujol, 853"

# Row 39
$ws.Cells.Item(39, 1).Value = "La Braciera Pizzaria - Pizza Napoletana"
$ws.Cells.Item(39, 3).Value = "4,8"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "4.166"
$ws.Cells.Item(39, 5).Value = "R. Conselheiro Saraiva, 664"

# Row 40
$ws.Cells.Item(40, 3).Value = "4,3"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "572"
$ws.Cells.Item(40, 5).Value = "R. Voluntários da Pátria, 3670"

# Row 41
$ws.Cells.Item(41, 3).Value = "4,8"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "45.025"
$ws.Cells.Item(41, 5).Value = "Av. Braz Leme, 201"

# Row 42
$ws.Cells.Item(42, 5).Value = "R. Augusto Tolle, 279"

# Row 43
$ws.Cells.Item(43, 1).Value = "Paleteria paulista - unidade Santana"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "99"
$ws.Cells.Item(43, 5).Value = "Rua Dr. Olavo Egídio, 798"

# Row 44
$ws.Cells.Item(44, 1).Value = "Latte di Fiori"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "416"
$ws.Cells.Item(44, 5).Value = "Av. Águas de São Pedro, 427"

# Row 45
$ws.Cells.Item(45, 3).Value = "4,6"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "69"
$ws.Cells.Item(45, 5).Value = "Av. Nova Cantareira, 346"

# Row 46
$ws.Cells.Item(46, 3).Value = "4,6"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "304"
$ws.Cells.Item(46, 5).Value = "R. Dom Henrique Mourão, 216"

# Row 47
$ws.Cells.Item(47, 1).Value = "Gelato Borelli Santana"
$ws.Cells.Item(47, 3).Value = "4,5"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "431"
$ws.Cells.Item(47, 5).Value = "R. Tupiguaes, 48 - Loja 02"

# Row 48
$ws.Cells.Item(48, 1).Value = "Bacio di Latte - Augusto Tolle"
$ws.Cells.Item(48, 3).Value = "4,4"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "143"
$ws.Cells.Item(48, 5).Value = "R. Augusto Tolle, 619"

# Row 49
$ws.Cells.Item(49, 1).Value = "Lambisk Sorvetes"
$ws.Cells.Item(49, 3).Value = "4,7"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "192"
$ws.Cells.Item(49, 5).Value = "Rua Conselheiro Moreira de Barros, 1002"

# Row 50
$ws.Cells.Item(50, 1).Value = "Frida & Mina Sorveteria"
$ws.Cells.Item(50, 3).Value = "4,6"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "6.987"
$ws.Cells.Item(50, 5).Value = "R. Artur de Azevedo, 1147"

# Row 51
$ws.Cells.Item(51, 1).Value = "OGGI SORVETES PARADA INGLESA"
$ws.Cells.Item(51, 3).Value = "4,6"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "378"
$ws.Cells.Item(51, 5).Value = "Av. Gen. Ataliba Leonel, 1976"

# Row 52
$ws.Cells.Item(52, 1).Value = "Dezato Gelato - Braz Leme"
$ws.Cells.Item(52, 4).NumberFormat = "@"
$ws.Cells.Item(52, 4).Value = "260"
$ws.Cells.Item(52, 5).Value = "Av. Braz Leme, 1200 - Loja 06"

# Row 53
$ws.Cells.Item(53, 1).Value = "Sorveteria do Centro"
$ws.Cells.Item(53, 3).Value = "4,4"
$ws.Cells.Item(53, 4).NumberFormat = "@"
$ws.Cells.Item(53, 4).Value = "2.875"
$ws.Cells.Item(53, 5).Value = "R. Epitácio Pessoa, 94"

# Row 54
$ws.Cells.Item(54, 3).Value = "4,7"
$ws.Cells.Item(54, 4).NumberFormat = "@"
$ws.Cells.Item(54, 4).Value = "31"
$ws.Cells.Item(54, 5).Value = "Av. Nova Cantareira, 335"

# Row 55
$ws.Cells.Item(55, 1).Value = "Candy Village Sorvete Infancia vintage"
$ws.Cells.Item(55, 3).Value = "4,9"
$ws.Cells.Item(55, 4).NumberFormat = "@"
$ws.Cells.Item(55, 4).Value = "33"
$ws.Cells.Item(55, 5).Value = "R. Jovita"

# Row 56
$ws.Cells.Item(56, 1).Value = "Sorveteria Lafer & Açaí"
$ws.Cells.Item(56, 3).Value = "4,6"
$ws.Cells.Item(56, 4).NumberFormat = "@"
$ws.Cells.Item(56, 4).Value = "982"
$ws.Cells.Item(56, 5).Value = "R. Maria Curupaiti, 1520"

# Row 57
$ws.Cells.Item(57, 1).Value = "Sou Ice Santa Teresinha - Sorvetes e Congelados"
$ws.Cells.Item(57, 4).NumberFormat = "@"
$ws.Cells.Item(57, 4).Value = "49"
$ws.Cells.Item(57, 5).Value = "R. Conselheiro Moreira de Barros, 1309 - Loja 01"

# Row 58
$ws.Cells.Item(58, 3).Value = "4,4"
$ws.Cells.Item(58, 4).NumberFormat = "@"
$ws.Cells.Item(58, 4).Value = "1.593"
$ws.Cells.Item(58, 5).Value = "Av. Braz Leme, 2378"

# Row 59
$ws.Cells.Item(59, 3).Value = "4,7"
$ws.Cells.Item(59, 4).NumberFormat = "@"
$ws.Cells.Item(59, 4).Value = "135"
$ws.Cells.Item(59, 5).Value = "Av. Conceição, 898"

# Row 60
$ws.Cells.Item(60, 3).Value = "4,9"
$ws.Cells.Item(60, 4).NumberFormat = "@"
$ws.Cells.Item(60, 4).Value = "87"
$ws.Cells.Item(60, 5).Value = "Av. Gen. Ataliba Leonel, 3837"

# Row 61
$ws.Cells.Item(61, 1).Value = "Bacio di Latte"
$ws.Cells.Item(61, 3).Value = "4,5"
$ws.Cells.Item(61, 4).NumberFormat = "@"
$ws.Cells.Item(61, 4).Value = "807"
$ws.Cells.Item(61, 5).Value = "Tv. Casalbuono, 120"

